$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.355422377586365
$ws.Range("B1").Value = 2.014503717422485
$ws.Range("C1").Value = 3.854816198348999
$ws.Range("D1").Value = 1.017550110816956
$ws.Range("E1").Value = 0.7481813430786133
